# Update the slide resolution (slide size) for the presentation and
# resize every slide's full-bleed picture to match the new canvas.
#
# Target size: cx=14630400 EMU (1152 pt), cy=8252460 EMU (649.8 pt)
#
# Note: PageSetup.SlideWidth/SlideHeight round-to-nearest-EMU, so the
# "natural" point values (1152 / 649.8) land exactly on the target EMU
# values. Shape.Width/Height truncate towards zero, so the height needs
# a tiny nudge (649.80002) to avoid landing one EMU short (8252459
# instead of 8252460).

$p = $ppt.ActivePresentation

$newWidthPt  = 1152
$newHeightPt = 649.80002

# 1) Resize the presentation canvas itself.
$p.PageSetup.SlideWidth  = 1152
$p.PageSetup.SlideHeight = 649.8

# 2) Resize every slide's full-bleed picture (offset stays at 0,0) so it
#    keeps covering the whole slide at the new resolution.
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        $shape.Left   = 0
        $shape.Top    = 0
        $shape.Width  = $newWidthPt
        $shape.Height = $newHeightPt
    }
}
